# Auto-generated edit script applying the diffed changes to before.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = 47438
$ws.Range("C15").Value = "SIG-3w Lilliput LED Torch &amp; Table Lamp"
$ws.Range("D15").Value = 401.81
$ws.Range("E15").Value = 480.05
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 803.62
$ws.Range("B16").Value = 59408
$ws.Range("C16").Value = "SIG-3W Lilliput LED Torch &amp; Table Lamp"
$ws.Range("D16").Value = 388.17
$ws.Range("E16").Value = 463.78
$ws.Range("F16").Value = 38
$ws.Range("G16").Value = 14750.46
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("F22").Value = 80
$ws.Range("G22").Value = 11647.2
$ws.Range("B25").Value = 79876.74000000001
$ws.Range("F42").Value = 4
$ws.Range("G42").Value = 10809.28
$ws.Range("B47").Value = 58083.24
$ws.Range("F61").Value = 4
$ws.Range("G61").Value = 380.08
$ws.Range("F63").Value = 29
$ws.Range("G63").Value = 2036.96
$ws.Range("F69").Value = 311
$ws.Range("G69").Value = 34887.98
$ws.Range("F74").Value = 7
$ws.Range("G74").Value = 159.18
$ws.Range("F79").Value = 20
$ws.Range("G79").Value = 1881.8
$ws.Range("F82").Value = 40
$ws.Range("G82").Value = 661.6
$ws.Range("F83").Value = 29
$ws.Range("G83").Value = 1802.93
$ws.Range("F92").Value = 70
$ws.Range("G92").Value = 3129
$ws.Range("B95").Value = 127100.76
$ws.Range("F126").Value = 73
$ws.Range("G126").Value = 3612.04
$ws.Range("F129").Value = 3
$ws.Range("G129").Value = 126.96
$ws.Range("B135").Value = 17724.57
$ws.Range("F143").Value = 2
$ws.Range("G143").Value = 86.26000000000001
$ws.Range("F151").Value = 1
$ws.Range("G151").Value = 33.01
$ws.Range("F159").Value = 1
$ws.Range("G159").Value = 54.41
$ws.Range("B168").Value = 31707.05
$ws.Range("F202").Value = 170
$ws.Range("G202").Value = 12192.4
$ws.Range("F203").Value = 48
$ws.Range("G203").Value = 4276.32
$ws.Range("F206").Value = 16
$ws.Range("G206").Value = 754.4
$ws.Range("F207").Value = 116
$ws.Range("G207").Value = 8618.799999999999
$ws.Range("B208").Value = 33468.64
$ws.Range("F210").Value = 153
$ws.Range("G210").Value = 17888.76
$ws.Range("F221").Value = 24
$ws.Range("G221").Value = 1069.92
$ws.Range("B222").Value = 51889.51
$ws.Range("F237").Value = 6
$ws.Range("G237").Value = 192.42
$ws.Range("F240").Value = 34
$ws.Range("G240").Value = 1062.5
$ws.Range("F241").Value = 60
$ws.Range("G241").Value = 17745
$ws.Range("F247").Value = 31
$ws.Range("G247").Value = 1167.15
$ws.Range("F250").Value = 4
$ws.Range("G250").Value = 187.48
$ws.Range("F257").Value = 0
$ws.Range("G257").Value = 0
$ws.Range("B258").Value = 42410.91
$ws.Range("F272").Value = 32
$ws.Range("G272").Value = 1386.24
$ws.Range("F277").Value = 61
$ws.Range("G277").Value = 2937.15
$ws.Range("F289").Value = 85
$ws.Range("G289").Value = 2245.7
$ws.Range("B290").Value = 70046.07000000001
$ws.Range("F305").Value = 38
$ws.Range("G305").Value = 1021.82
$ws.Range("B307").Value = 7086.74
$ws.Range("F309").Value = 55
$ws.Range("G309").Value = 3657.5
$ws.Range("F311").Value = 34
$ws.Range("G311").Value = 1696.26
$ws.Range("F320").Value = 49
$ws.Range("G320").Value = 9027.27
$ws.Range("F330").Value = 17
$ws.Range("G330").Value = 805.8
$ws.Range("F337").Value = 112
$ws.Range("G337").Value = 8060.64
$ws.Range("F341").Value = 22
$ws.Range("G341").Value = 1097.58
$ws.Range("B347").Value = 120865.08
$ws.Range("F382").Value = 9
$ws.Range("G382").Value = 372.78
$ws.Range("F386").Value = 113
$ws.Range("G386").Value = 1066.72
$ws.Range("F387").Value = 70
$ws.Range("G387").Value = 1541.4
$ws.Range("F392").Value = 16
$ws.Range("G392").Value = 1554.72
$ws.Range("B393").Value = 8630.65
$ws.Range("F423").Value = 16
$ws.Range("G423").Value = 2842.56
$ws.Range("B425").Value = 2911.09
$ws.Range("F435").Value = 28
$ws.Range("G435").Value = 1430.24
$ws.Range("F436").Value = 115
$ws.Range("G436").Value = 5776.45
$ws.Range("F439").Value = 31
$ws.Range("G439").Value = 1701.28
$ws.Range("B455").Value = 43026.83
$ws.Range("F463").Value = 55
$ws.Range("G463").Value = 1821.05
$ws.Range("B470").Value = 10161.24
$ws.Range("F596").Value = 73
$ws.Range("G596").Value = 2929.49
$ws.Range("F599").Value = 101
$ws.Range("G599").Value = 3931.93
$ws.Range("F600").Value = 117
$ws.Range("G600").Value = 4625.01
$ws.Range("B601").Value = 20556.18
$ws.Range("B607").Value = 1568100.16
$ws.Range("B608").Value = 1568100.16
